$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column D (4) so the stored OOXML width becomes 53
# (Excel's ColumnWidth property is offset by 5/6 of a character from the
# width that ends up written into <col width="...">.)
$ws.Columns.Item(4).ColumnWidth = 52.1666666666667

# Append a new run-log row (row 3)
$ws.Range("A3").Value = "2025-08-12 06:51:35 UTC"
$ws.Range("B3").Value = "2025-08-12 12:21:35 IST"
$ws.Range("C3").Value = "SKIPPED"
$ws.Range("D3").Value = "No change in PDF. Skipping download & Excel update."
$ws.Range("E3").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = ""

# Match the formatting used by the previous data row (row 2) so the new
# row reuses the same cell style (center/center, no border/bold)
$ws.Range("A2:H2").Copy()
$ws.Range("A3:H3").PasteSpecial(-4122)
